$wb = $excel.ActiveWorkbook
$wsDBD = $wb.Worksheets.Item("DBD")
$wsDBS = $wb.Worksheets.Item("DBS")

# 1) DBD sheet: field #6 "Entdy" Chinese name changes from "交易日" (transaction date)
#    to "會計日" (accounting date)
$wsDBD.Range("C14").Value = "會計日"

# 2) DBS sheet: add a brand-new row describing a "findCreatDate" lookup function
#    (values entered before the B3 update below so new shared strings line up the
#    same way Excel itself ordered them)
$wsDBS.Range("A4").Value = "findCreatDate"
$wsDBS.Range("C4").Value = "CreateDate"

# 3) DBS sheet: findEntdy's filter text gains an extra "AND SupNo %" clause, and the
#    new row's filter text is filled in
$wsDBS.Range("B3").Value = "Entdy >= ,AND Entdy <=,AND SupNo %"
$wsDBS.Range("B4").Value = "CreateDate>=, AND CreateDate<= ,AND SupNo %"

# Match the formatting of the new B4 cell (left/center aligned, no wrap) like the rest of column B
$wsDBS.Range("B4").WrapText = $false
$wsDBS.Range("B4").HorizontalAlignment = -4131
$wsDBS.Range("B4").VerticalAlignment = -4108

# Restore the last-selected cells recorded in each sheet view
$wsDBD.Range("C17").Select()
$wsDBS.Range("B5").Select()
